# "Generate Report for Archive"
#
# 1. The localization status text changes from "Ready for handoff" to
#    "In Translation" everywhere it is used:
#      - Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3
#      - zh-cn sheet:     column C (Status), rows 2-3
#      - de-de sheet:     column C (Status), rows 2-3
# 2. The now-narrower status column is resized on all three sheets
#    (Overview columns E & F, and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "In Translation"
$newColWidth = 12.5

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $statusNew
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = $statusNew
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = $statusNew
$dede.Columns.Item(3).ColumnWidth = $newColWidth
